# Added general monitoring to Section 2.
# - Move the "Symbol Type" rectangle left to make room for a new shape.
# - Duplicate it to create a new "Objects" rectangle box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the existing "Symbol Type" rectangle (cNvPr id="87", name="Rectangle 9")
# and shift it left from x=6400800 EMU (504pt) to x=5486400 EMU (432pt).
$symbolShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 87) {
        $symbolShape = $candidate
        break
    }
}

$symbolShape.Left = 432

# Duplicate the "Symbol Type" shape (carries over its style/formatting) to
# build the new "Objects" rectangle, then move/resize it into place and
# replace its text.
$newShape = $symbolShape.Duplicate()
$newShape.Left = 558
$newShape.Top = 91.70275590551181
$newShape.Width = 86.98377952755905
$newShape.Height = 58.29724409448819
$newShape.TextFrame.TextRange.Text = "Objects"
